$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update status text from "In Translation" to "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Update timestamps
$overview.Range("G2").Value = "2016-08-25 09:03:15"
$dede.Range("H2").Value = "2016-08-25 09:03:15"
$zhcn.Range("H2").Value = "2016-08-25 09:03:05"

# Adjust affected column widths to reflect the new, longer status text
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3
